# Attendance update (9th Feb) - mark Day 5 (column K) attendance for each
# participant row 7-82 on the SM-II (F) sheet: "P" for present, "A" for the
# participants who were absent that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K currently has no value (just inherited blank-cell styling from
# columns L..AA). Copy the formatting Excel already uses for the filled-in
# attendance columns (column J, "Day 4") down into column K so the new
# entries pick up the same look-and-feel the other marked days have.
$ws.Range("J7:J82").Copy() | Out-Null
$ws.Range("K7:K82").PasteSpecial(-4122) | Out-Null

# Roll numbers (by row) marked absent ("A") on Day 5; everyone else attended
# ("P") that day.
$absentRows = @(10, 17, 32, 40, 42, 46, 47, 55, 57, 64, 66, 67, 72, 81)

for ($row = 7; $row -le 82; $row++) {
    if ($absentRows -contains $row) {
        $ws.Cells.Item($row, 11).Value = "A"
    } else {
        $ws.Cells.Item($row, 11).Value = "P"
    }
}
